# Auto-generated script to apply cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to
# text format first, otherwise Excel auto-converts them to numeric cells
# (losing the original inlineStr/text representation used in the sheet).
$textForceCells = @(
    'D5',
    'D6',
    'D7',
    'D10',
    'D11',
    'D12',
    'D14',
    'D19',
    'D20',
    'D21',
    'D22',
    'D23',
    'D24',
    'D25',
    'D26',
    'D28',
    'D29',
    'D30',
    'D31',
    'D33',
    'D34',
    'D35',
    'D36',
    'D37',
    'D38',
    'D39',
    'D42',
    'D44',
    'D45',
    'D50',
)
foreach ($c in $textForceCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply the updated values
$ws.Range('D2').Value = '64.507.36'
$ws.Range('E2').Value = '  +5.56%  '
$ws.Range('D3').Value = '3.083.03'
$ws.Range('E3').Value = '  +3.78%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').Value = '554.70'
$ws.Range('E5').Value = '  +1.93%  '
$ws.Range('D6').Value = '143.44'
$ws.Range('E6').Value = '  +10.53%  '
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('D8').Value = '3.081.81'
$ws.Range('E8').Value = '  +3.79%  '
$ws.Range('E9').Value = '  +1.17%  '
$ws.Range('D10').Value = '6.96'
$ws.Range('E10').Value = '  +18.89%  '
$ws.Range('D11').Value = '0.152'
$ws.Range('E11').Value = '  +6.05%  '
$ws.Range('D12').Value = '0.456'
$ws.Range('E12').Value = '  +3.72%  '
$ws.Range('E13').Value = '  +4.62%  '
$ws.Range('D14').Value = '35.26'
$ws.Range('E14').Value = '  +5.77%  '
$ws.Range('D15').Value = '3.575.37'
$ws.Range('E15').Value = '  +3.37%  '
$ws.Range('D16').Value = '64.474.30'
$ws.Range('E16').Value = '  +5.48%  '
$ws.Range('D17').Value = '3.078.87'
$ws.Range('E17').Value = '  +3.50%  '
$ws.Range('E18').Value = '  -1.07%  '
$ws.Range('D19').Value = '6.72'
$ws.Range('E19').Value = '  +3.40%  '
$ws.Range('D20').Value = '480.11'
$ws.Range('E20').Value = '  +2.50%  '
$ws.Range('D21').Value = '13.66'
$ws.Range('E21').Value = '  +5.66%  '
$ws.Range('D22').Value = '0.671'
$ws.Range('E22').Value = '  +2.15%  '
$ws.Range('D23').Value = '7.47'
$ws.Range('E23').Value = '  +8.95%  '
$ws.Range('D24').Value = '13.14'
$ws.Range('E24').Value = '  +11.29%  '
$ws.Range('D25').Value = '80.63'
$ws.Range('E25').Value = '  +1.79%  '
$ws.Range('D26').Value = '1.01'
$ws.Range('E26').Value = '  +0.55%  '
$ws.Range('E27').Value = '  +4.47%  '
$ws.Range('D28').Value = '7.91'
$ws.Range('E28').Value = '  +5.29%  '
$ws.Range('D29').Value = '2.06'
$ws.Range('E29').Value = '  +10.52%  '
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').Value = '  -0.15%  '
$ws.Range('D31').Value = '26.13'
$ws.Range('E31').Value = '  +3.83%  '
$ws.Range('E32').Value = '  +3.96%  '
$ws.Range('D33').Value = '2.43'
$ws.Range('E33').Value = '  +7.78%  '
$ws.Range('D34').Value = '5.71'
$ws.Range('E34').Value = '  +6.10%  '
$ws.Range('D35').Value = '55.23'
$ws.Range('E35').Value = '  +1.96%  '
$ws.Range('D36').Value = '6.08'
$ws.Range('E36').Value = '  +5.27%  '
$ws.Range('D37').Value = '464.06'
$ws.Range('E37').Value = '  +5.55%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').Value = '0.0827'
$ws.Range('E38').Value = '  +6.22%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = '0.0406'
$ws.Range('E39').Value = '  +8.77%  '
$ws.Range('D40').Value = '3.010.73'
$ws.Range('E40').Value = '  -3.14%  '
$ws.Range('E41').Value = '  +1.53%  '
$ws.Range('D42').Value = '8.26'
$ws.Range('E42').Value = '  +3.62%  '
$ws.Range('E43').Value = '  +16.43%  '
$ws.Range('D44').Value = '27.78'
$ws.Range('E44').Value = '  +10.55%  '
$ws.Range('D45').Value = '0.257'
$ws.Range('E45').Value = '  +8.40%  '
$ws.Range('E47').Value = '  +7.80%  '
$ws.Range('E48').Value = '  +4.59%  '
$ws.Range('D49').Value = '0.0₃0514'
$ws.Range('E49').Value = '  +8.53%  '
$ws.Range('D50').Value = '116.85'
$ws.Range('E50').Value = '  +3.03%  '
$ws.Range('E51').Value = '  +5.10%  '

# Restore default cell style for the cells we forced to text format,
# so no stray number-format style lingers on them.
foreach ($c in $textForceCells) {
    $ws.Range($c).Style = "Normal"
}
